# Auto-generated Excel COM-interop script to apply the Faerie_Profits data refresh diff.
# Updates computed marketboard/profit values across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 8).Value = 384.15  # H2: 380.0476 -> 384.15
$ws.Cells.Item(2, 10).Value = 964.4  # J2: 853.3333 -> 964.4
$ws.Cells.Item(2, 12).Value = 964.4  # L2: 853.3333 -> 964.4
$ws.Cells.Item(2, 14).Value = -1190.4  # N2: -1079.3333 -> -1190.4
$ws.Cells.Item(17, 8).Value = 142684.52  # H17: 7929.7427 -> 142684.52
$ws.Cells.Item(17, 10).Value = 156541.47  # J17: 251.36667 -> 156541.47
$ws.Cells.Item(17, 12).Value = 469624.41  # L17: 754.10001 -> 469624.41
$ws.Cells.Item(17, 14).Value = -469960.41  # N17: -1090.10001 -> -469960.41
$ws.Cells.Item(18, 8).Value = 3548.1875  # H18: 3198.4443 -> 3548.1875
$ws.Cells.Item(18, 9).Value = 2776.6155  # I18: 2599.6428 -> 2776.6155
$ws.Cells.Item(18, 10).Value = 6891.6665  # J18: 5294.25 -> 6891.6665
$ws.Cells.Item(18, 11).Value = 2776.6155  # K18: 2599.6428 -> 2776.6155
$ws.Cells.Item(18, 12).Value = 6891.6665  # L18: 5294.25 -> 6891.6665
$ws.Cells.Item(18, 13).Value = -2492.6155  # M18: -2315.6428 -> -2492.6155
$ws.Cells.Item(18, 14).Value = -7459.6665  # N18: -5862.25 -> -7459.6665
$ws.Cells.Item(26, 8).Value = 0  # H26: 49999 -> 0
$ws.Cells.Item(26, 10).Value = 0  # J26: 49999 -> 0
$ws.Cells.Item(26, 12).Value = 0  # L26: 49999 -> 0
$ws.Cells.Item(26, 14).ClearContents()  # N26: -50687 -> (removed)
$ws.Cells.Item(28, 8).Value = 1167.4138  # H28: 1154.2593 -> 1167.4138
$ws.Cells.Item(28, 9).Value = 833.0952  # I28: 779.2105 -> 833.0952
$ws.Cells.Item(28, 11).Value = 833.0952  # K28: 779.2105 -> 833.0952
$ws.Cells.Item(28, 13).Value = -348.0952  # M28: -294.2105 -> -348.0952
$ws.Cells.Item(29, 8).Value = 3800.6667  # H29: 4850.5 -> 3800.6667
$ws.Cells.Item(29, 10).Value = 7999  # J29: 7999.5 -> 7999
$ws.Cells.Item(29, 12).Value = 23997  # L29: 23998.5 -> 23997
$ws.Cells.Item(29, 14).Value = -24559  # N29: -24560.5 -> -24559
$ws.Cells.Item(32, 8).Value = 71436380  # H32: 50005696 -> 71436380
$ws.Cells.Item(32, 9).Value = 142867820  # I32: 125009430 -> 142867820
$ws.Cells.Item(32, 10).Value = 4948.5713  # J32: 3207.9167 -> 4948.5713
$ws.Cells.Item(32, 11).Value = 142867820  # K32: 125009430 -> 142867820
$ws.Cells.Item(32, 12).Value = 4948.5713  # L32: 3207.9167 -> 4948.5713
$ws.Cells.Item(32, 13).Value = -142867494  # M32: -125009104 -> -142867494
$ws.Cells.Item(32, 14).Value = -5600.5713  # N32: -3859.9167 -> -5600.5713
$ws.Cells.Item(33, 8).Value = 380.52942  # H33: 402.875 -> 380.52942
$ws.Cells.Item(33, 9).Value = 348.08334  # I33: 377.63635 -> 348.08334
$ws.Cells.Item(33, 11).Value = 348.08334  # K33: 377.63635 -> 348.08334
$ws.Cells.Item(33, 13).Value = -119.08334  # M33: -148.63635 -> -119.08334
$ws.Cells.Item(38, 8).Value = 1918.8  # H38: 2028.1666 -> 1918.8
$ws.Cells.Item(38, 9).Value = 773.5  # I38: 704.2222 -> 773.5
$ws.Cells.Item(38, 10).Value = 6500  # J38: 6000 -> 6500
$ws.Cells.Item(38, 11).Value = 2320.5  # K38: 2112.6666 -> 2320.5
$ws.Cells.Item(38, 12).Value = 19500  # L38: 18000 -> 19500
$ws.Cells.Item(38, 13).Value = -1948.5  # M38: -1740.6666 -> -1948.5
$ws.Cells.Item(38, 14).Value = -20244  # N38: -18744 -> -20244
$ws.Cells.Item(51, 8).Value = 71431624  # H51: 55558260 -> 71431624
$ws.Cells.Item(51, 9).Value = 166668340  # I51: 125001624 -> 166668340
$ws.Cells.Item(51, 10).Value = 4086.375  # J51: 3564.8 -> 4086.375
$ws.Cells.Item(51, 11).Value = 166668340  # K51: 125001624 -> 166668340
$ws.Cells.Item(51, 12).Value = 4086.375  # L51: 3564.8 -> 4086.375
$ws.Cells.Item(51, 13).Value = -166667856  # M51: -125001140 -> -166667856
$ws.Cells.Item(51, 14).Value = -5054.375  # N51: -4532.8 -> -5054.375
$ws.Cells.Item(58, 8).Value = 253  # H58: 335 -> 253
$ws.Cells.Item(58, 9).Value = 253  # I58: 225 -> 253
$ws.Cells.Item(58, 10).Value = 0  # J58: 995 -> 0
$ws.Cells.Item(58, 11).Value = 759  # K58: 675 -> 759
$ws.Cells.Item(58, 12).Value = 0  # L58: 2985 -> 0
$ws.Cells.Item(58, 13).Value = -609  # M58: -525 -> -609
$ws.Cells.Item(58, 14).ClearContents()  # N58: -3285 -> (removed)
$ws.Cells.Item(62, 8).Value = 46779.68  # H62: 43537.11 -> 46779.68
$ws.Cells.Item(62, 9).Value = 62516.168  # I62: 59331.633 -> 62516.168
$ws.Cells.Item(62, 10).Value = 6314.4287  # J62: 6025.125 -> 6314.4287
$ws.Cells.Item(62, 11).Value = 62516.168  # K62: 59331.633 -> 62516.168
$ws.Cells.Item(62, 12).Value = 6314.4287  # L62: 6025.125 -> 6314.4287
$ws.Cells.Item(62, 13).Value = -61892.168  # M62: -58707.633 -> -61892.168
$ws.Cells.Item(62, 14).Value = -7562.4287  # N62: -7273.125 -> -7562.4287
$ws.Cells.Item(65, 8).Value = 46779.68  # H65: 43537.11 -> 46779.68
$ws.Cells.Item(65, 9).Value = 62516.168  # I65: 59331.633 -> 62516.168
$ws.Cells.Item(65, 10).Value = 6314.4287  # J65: 6025.125 -> 6314.4287
$ws.Cells.Item(65, 11).Value = 312580.84  # K65: 296658.165 -> 312580.84
$ws.Cells.Item(65, 12).Value = 31572.1435  # L65: 30125.625 -> 31572.1435
$ws.Cells.Item(65, 13).Value = -309460.84  # M65: -293538.165 -> -309460.84
$ws.Cells.Item(65, 14).Value = -37812.14350000001  # N65: -36365.625 -> -37812.14350000001
$ws.Cells.Item(70, 8).Value = 51159.184  # H70: 53528.668 -> 51159.184
$ws.Cells.Item(70, 10).Value = 7991.8335  # J70: 8591.091 -> 7991.8335
$ws.Cells.Item(70, 12).Value = 23975.5005  # L70: 25773.273 -> 23975.5005
$ws.Cells.Item(70, 14).Value = -24515.5005  # N70: -26313.273 -> -24515.5005
$ws.Cells.Item(73, 8).Value = 51159.184  # H73: 53528.668 -> 51159.184
$ws.Cells.Item(73, 10).Value = 7991.8335  # J73: 8591.091 -> 7991.8335
$ws.Cells.Item(73, 12).Value = 23975.5005  # L73: 25773.273 -> 23975.5005
$ws.Cells.Item(73, 14).Value = -25847.5005  # N73: -27645.273 -> -25847.5005
$ws.Cells.Item(76, 8).Value = 3172.0908  # H76: 2965 -> 3172.0908
$ws.Cells.Item(76, 9).Value = 3172.0908  # I76: 2965 -> 3172.0908
$ws.Cells.Item(76, 11).Value = 3172.0908  # K76: 2965 -> 3172.0908
$ws.Cells.Item(76, 13).Value = -2857.0908  # M76: -2650 -> -2857.0908
$ws.Cells.Item(79, 8).Value = 3172.0908  # H79: 2965 -> 3172.0908
$ws.Cells.Item(79, 9).Value = 3172.0908  # I79: 2965 -> 3172.0908
$ws.Cells.Item(79, 11).Value = 3172.0908  # K79: 2965 -> 3172.0908
$ws.Cells.Item(79, 13).Value = -2080.0908  # M79: -1873 -> -2080.0908
$ws.Cells.Item(86, 8).Value = 3380.1333  # H86: 3348.484 -> 3380.1333
$ws.Cells.Item(86, 10).Value = 3580.5334  # J86: 3506.6875 -> 3580.5334
$ws.Cells.Item(86, 12).Value = 3580.5334  # L86: 3506.6875 -> 3580.5334
$ws.Cells.Item(86, 14).Value = -5826.5334  # N86: -5752.6875 -> -5826.5334
$ws.Cells.Item(88, 8).Value = 9133.9  # H88: 9949.111000000001 -> 9133.9
$ws.Cells.Item(88, 9).Value = 8590.333000000001  # I88: 8923.833000000001 -> 8590.333000000001
$ws.Cells.Item(88, 10).Value = 9949.25  # J88: 11999.667 -> 9949.25
$ws.Cells.Item(88, 11).Value = 8590.333000000001  # K88: 8923.833000000001 -> 8590.333000000001
$ws.Cells.Item(88, 12).Value = 9949.25  # L88: 11999.667 -> 9949.25
$ws.Cells.Item(88, 13).Value = -8184.333000000001  # M88: -8517.833000000001 -> -8184.333000000001
$ws.Cells.Item(88, 14).Value = -10761.25  # N88: -12811.667 -> -10761.25
$ws.Cells.Item(89, 8).Value = 3380.1333  # H89: 3348.484 -> 3380.1333
$ws.Cells.Item(89, 10).Value = 3580.5334  # J89: 3506.6875 -> 3580.5334
$ws.Cells.Item(89, 12).Value = 17902.667  # L89: 17533.4375 -> 17902.667
$ws.Cells.Item(89, 14).Value = -29134.667  # N89: -28765.4375 -> -29134.667
$ws.Cells.Item(91, 8).Value = 9133.9  # H91: 9949.111000000001 -> 9133.9
$ws.Cells.Item(91, 9).Value = 8590.333000000001  # I91: 8923.833000000001 -> 8590.333000000001
$ws.Cells.Item(91, 10).Value = 9949.25  # J91: 11999.667 -> 9949.25
$ws.Cells.Item(91, 11).Value = 8590.333000000001  # K91: 8923.833000000001 -> 8590.333000000001
$ws.Cells.Item(91, 12).Value = 9949.25  # L91: 11999.667 -> 9949.25
$ws.Cells.Item(91, 13).Value = -7186.333000000001  # M91: -7519.833000000001 -> -7186.333000000001
$ws.Cells.Item(91, 14).Value = -12757.25  # N91: -14807.667 -> -12757.25
$ws.Cells.Item(98, 8).Value = 5954761.5  # H98: 5954770 -> 5954761.5
$ws.Cells.Item(98, 9).Value = 7144623.5  # I98: 7144633.5 -> 7144623.5
$ws.Cells.Item(98, 11).Value = 7144623.5  # K98: 7144633.5 -> 7144623.5
$ws.Cells.Item(98, 13).Value = -7143125.5  # M98: -7143135.5 -> -7143125.5
$ws.Cells.Item(103, 8).Value = 1862  # H103: 1564.25 -> 1862
$ws.Cells.Item(103, 9).Value = 816  # I103: 921.44446 -> 816
$ws.Cells.Item(103, 10).Value = 5000  # J103: 3492.6667 -> 5000
$ws.Cells.Item(103, 11).Value = 2448  # K103: 2764.33338 -> 2448
$ws.Cells.Item(103, 12).Value = 15000  # L103: 10478.0001 -> 15000
$ws.Cells.Item(103, 13).Value = -1862  # M103: -2178.33338 -> -1862
$ws.Cells.Item(103, 14).Value = -16172  # N103: -11650.0001 -> -16172
$ws.Cells.Item(106, 8).Value = 51675.61  # H106: 56073.57 -> 51675.61
$ws.Cells.Item(106, 9).Value = 56427  # I106: 62085.89 -> 56427
$ws.Cells.Item(106, 11).Value = 56427  # K106: 62085.89 -> 56427
$ws.Cells.Item(106, 13).Value = -55796  # M106: -61454.89 -> -55796
$ws.Cells.Item(107, 8).Value = 984.7714  # H107: 807.96875 -> 984.7714
$ws.Cells.Item(107, 9).Value = 438.625  # I107: 448.9565 -> 438.625
$ws.Cells.Item(107, 10).Value = 2176.3635  # J107: 1725.4445 -> 2176.3635
$ws.Cells.Item(107, 11).Value = 438.625  # K107: 448.9565 -> 438.625
$ws.Cells.Item(107, 12).Value = 2176.3635  # L107: 1725.4445 -> 2176.3635
$ws.Cells.Item(107, 13).Value = 1481.375  # M107: 1471.0435 -> 1481.375
$ws.Cells.Item(107, 14).Value = -6016.363499999999  # N107: -5565.4445 -> -6016.363499999999
$ws.Cells.Item(112, 8).Value = 627692.8  # H112: 558226.9 -> 627692.8
$ws.Cells.Item(112, 10).Value = 835924  # J112: 716863.4 -> 835924
$ws.Cells.Item(112, 12).Value = 2507772  # L112: 2150590.2 -> 2507772
$ws.Cells.Item(112, 14).Value = -2509988  # N112: -2152806.2 -> -2509988
$ws.Cells.Item(113, 8).Value = 3934  # H113: 3791.4348 -> 3934
$ws.Cells.Item(113, 9).Value = 2476.5833  # I113: 2450.5715 -> 2476.5833
$ws.Cells.Item(113, 11).Value = 2476.5833  # K113: 2450.5715 -> 2476.5833
$ws.Cells.Item(113, 13).Value = 777.4167000000002  # M113: 803.4285 -> 777.4167000000002
$ws.Cells.Item(122, 8).Value = 5954761.5  # H122: 5954770 -> 5954761.5
$ws.Cells.Item(122, 9).Value = 7144623.5  # I122: 7144633.5 -> 7144623.5
$ws.Cells.Item(122, 11).Value = 21433870.5  # K122: 21433900.5 -> 21433870.5
$ws.Cells.Item(122, 13).Value = -21431420.5  # M122: -21431450.5 -> -21431420.5
$ws.Cells.Item(132, 8).Value = 3188.449  # H132: 3677.1904 -> 3188.449
$ws.Cells.Item(132, 9).Value = 2984.0417  # I132: 3449.805 -> 2984.0417
$ws.Cells.Item(132, 11).Value = 8952.125100000001  # K132: 10349.415 -> 8952.125100000001
$ws.Cells.Item(132, 13).Value = -6422.125100000001  # M132: -7819.414999999999 -> -6422.125100000001
$ws.Cells.Item(135, 8).Value = 4957.48  # H135: 5227.7393 -> 4957.48
$ws.Cells.Item(135, 9).Value = 4414.0415  # I135: 4647.1816 -> 4414.0415
$ws.Cells.Item(135, 11).Value = 39726.3735  # K135: 41824.6344 -> 39726.3735
$ws.Cells.Item(135, 13).Value = -37191.3735  # M135: -39289.6344 -> -37191.3735
$ws.Cells.Item(137, 8).Value = 3548.725  # H137: 3709.5 -> 3548.725
$ws.Cells.Item(137, 9).Value = 4146.143  # I137: 4637.6113 -> 4146.143
$ws.Cells.Item(137, 10).Value = 2888.4211  # J137: 2874.2 -> 2888.4211
$ws.Cells.Item(137, 11).Value = 12438.429  # K137: 13912.8339 -> 12438.429
$ws.Cells.Item(137, 12).Value = 8665.263300000001  # L137: 8622.599999999999 -> 8665.263300000001
$ws.Cells.Item(137, 13).Value = -9888.429  # M137: -11362.8339 -> -9888.429
$ws.Cells.Item(137, 14).Value = -13765.2633  # N137: -13722.6 -> -13765.2633
$ws.Cells.Item(138, 8).Value = 913816.75  # H138: 719355.6 -> 913816.75
$ws.Cells.Item(138, 10).Value = 1672999.6  # J138: 1117443.6 -> 1672999.6
$ws.Cells.Item(138, 12).Value = 5018998.800000001  # L138: 3352330.8 -> 5018998.800000001
$ws.Cells.Item(138, 14).Value = -5029278.800000001  # N138: -3362610.8 -> -5029278.800000001
$ws.Cells.Item(141, 8).Value = 2429.8333  # H141: 2041.2222 -> 2429.8333
$ws.Cells.Item(141, 9).Value = 1914.8125  # I141: 1978.6471 -> 1914.8125
$ws.Cells.Item(141, 10).Value = 6550  # J141: 3105 -> 6550
$ws.Cells.Item(141, 11).Value = 5744.4375  # K141: 5935.9413 -> 5744.4375
$ws.Cells.Item(141, 12).Value = 19650  # L141: 9315 -> 19650
$ws.Cells.Item(141, 13).Value = -564.4375  # M141: -755.9412999999995 -> -564.4375
$ws.Cells.Item(141, 14).Value = -30010  # N141: -19675 -> -30010

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 1399.4844  # H32: 1446.5238 -> 1399.4844
$ws.Cells.Item(32, 9).Value = 1399.4678  # I32: 1448.0492 -> 1399.4678
$ws.Cells.Item(32, 11).Value = 1399.4678  # K32: 1448.0492 -> 1399.4678
$ws.Cells.Item(32, 13).Value = -1112.4678  # M32: -1161.0492 -> -1112.4678
$ws.Cells.Item(45, 8).Value = 6171.696  # H45: 5885 -> 6171.696
$ws.Cells.Item(45, 9).Value = 6015.7896  # I45: 5689.3335 -> 6015.7896
$ws.Cells.Item(45, 11).Value = 6015.7896  # K45: 5689.3335 -> 6015.7896
$ws.Cells.Item(45, 13).Value = -5638.7896  # M45: -5312.3335 -> -5638.7896
$ws.Cells.Item(74, 8).Value = 4403.273  # H74: 4669.05 -> 4403.273
$ws.Cells.Item(74, 9).Value = 4527.5386  # I74: 5033.364 -> 4527.5386
$ws.Cells.Item(74, 11).Value = 4527.5386  # K74: 5033.364 -> 4527.5386
$ws.Cells.Item(74, 13).Value = -3653.5386  # M74: -4159.364 -> -3653.5386
$ws.Cells.Item(77, 8).Value = 4403.273  # H77: 4669.05 -> 4403.273
$ws.Cells.Item(77, 9).Value = 4527.5386  # I77: 5033.364 -> 4527.5386
$ws.Cells.Item(77, 11).Value = 22637.693  # K77: 25166.82 -> 22637.693
$ws.Cells.Item(77, 13).Value = -18269.693  # M77: -20798.82 -> -18269.693
$ws.Cells.Item(88, 8).Value = 2663.25  # H88: 2590.5 -> 2663.25
$ws.Cells.Item(88, 9).Value = 2161.2  # I88: 2181 -> 2161.2
$ws.Cells.Item(88, 10).Value = 3500  # J88: 3000 -> 3500
$ws.Cells.Item(88, 11).Value = 2161.2  # K88: 2181 -> 2161.2
$ws.Cells.Item(88, 12).Value = 3500  # L88: 3000 -> 3500
$ws.Cells.Item(88, 13).Value = -1755.2  # M88: -1775 -> -1755.2
$ws.Cells.Item(88, 14).Value = -4312  # N88: -3812 -> -4312
$ws.Cells.Item(91, 8).Value = 2663.25  # H91: 2590.5 -> 2663.25
$ws.Cells.Item(91, 9).Value = 2161.2  # I91: 2181 -> 2161.2
$ws.Cells.Item(91, 10).Value = 3500  # J91: 3000 -> 3500
$ws.Cells.Item(91, 11).Value = 2161.2  # K91: 2181 -> 2161.2
$ws.Cells.Item(91, 12).Value = 3500  # L91: 3000 -> 3500
$ws.Cells.Item(91, 13).Value = -757.1999999999998  # M91: -777 -> -757.1999999999998
$ws.Cells.Item(91, 14).Value = -6308  # N91: -5808 -> -6308
$ws.Cells.Item(97, 8).Value = 2298.8857  # H97: 2389.1177 -> 2298.8857
$ws.Cells.Item(97, 9).Value = 1071.1613  # I97: 1123.3103 -> 1071.1613
$ws.Cells.Item(97, 10).Value = 11813.75  # J97: 9730.799999999999 -> 11813.75
$ws.Cells.Item(97, 11).Value = 1071.1613  # K97: 1123.3103 -> 1071.1613
$ws.Cells.Item(97, 12).Value = 11813.75  # L97: 9730.799999999999 -> 11813.75
$ws.Cells.Item(97, 13).Value = -575.1613  # M97: -627.3103000000001 -> -575.1613
$ws.Cells.Item(97, 14).Value = -12805.75  # N97: -10722.8 -> -12805.75
$ws.Cells.Item(102, 8).Value = 15386730  # H102: 16668803 -> 15386730
$ws.Cells.Item(102, 9).Value = 2608.625  # I102: 2713.4285 -> 2608.625
$ws.Cells.Item(102, 10).Value = 40001324  # J102: 40001330 -> 40001324
$ws.Cells.Item(102, 11).Value = 2608.625  # K102: 2713.4285 -> 2608.625
$ws.Cells.Item(102, 12).Value = 40001324  # L102: 40001330 -> 40001324
$ws.Cells.Item(102, 13).Value = -986.625  # M102: -1091.4285 -> -986.625
$ws.Cells.Item(102, 14).Value = -40004568  # N102: -40004574 -> -40004568
$ws.Cells.Item(110, 8).Value = 1419.8334  # H110: 1435.9166 -> 1419.8334
$ws.Cells.Item(110, 9).Value = 1331.7727  # I110: 1366.9546 -> 1331.7727
$ws.Cells.Item(110, 10).Value = 2388.5  # J110: 2194.5 -> 2388.5
$ws.Cells.Item(110, 11).Value = 1331.7727  # K110: 1366.9546 -> 1331.7727
$ws.Cells.Item(110, 12).Value = 2388.5  # L110: 2194.5 -> 2388.5
$ws.Cells.Item(110, 13).Value = 713.2273  # M110: 678.0454 -> 713.2273
$ws.Cells.Item(110, 14).Value = -6478.5  # N110: -6284.5 -> -6478.5
$ws.Cells.Item(122, 8).Value = 7493  # H122: 2103.3333 -> 7493
$ws.Cells.Item(122, 9).Value = 7493  # I122: 2103.3333 -> 7493
$ws.Cells.Item(122, 11).Value = 22479  # K122: 6309.999899999999 -> 22479
$ws.Cells.Item(122, 13).Value = -20029  # M122: -3859.999899999999 -> -20029
$ws.Cells.Item(124, 8).Value = 39950  # H124: 45000 -> 39950
$ws.Cells.Item(124, 10).Value = 39950  # J124: 45000 -> 39950
$ws.Cells.Item(124, 12).Value = 39950  # L124: 45000 -> 39950
$ws.Cells.Item(124, 14).Value = -49770  # N124: -54820 -> -49770
$ws.Cells.Item(125, 8).Value = 86717.164  # H125: 86797.60000000001 -> 86717.164
$ws.Cells.Item(125, 10).Value = 86717.164  # J125: 86797.60000000001 -> 86717.164
$ws.Cells.Item(125, 12).Value = 86717.164  # L125: 86797.60000000001 -> 86717.164
$ws.Cells.Item(125, 14).Value = -96557.164  # N125: -96637.60000000001 -> -96557.164
$ws.Cells.Item(132, 8).Value = 3438.9773  # H132: 3690.9534 -> 3438.9773
$ws.Cells.Item(132, 9).Value = 2064.8928  # I132: 2415.2964 -> 2064.8928
$ws.Cells.Item(132, 11).Value = 6194.678400000001  # K132: 7245.889200000001 -> 6194.678400000001
$ws.Cells.Item(132, 13).Value = -3664.678400000001  # M132: -4715.889200000001 -> -3664.678400000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(22, 8).Value = 242.42857  # H22: 227.25 -> 242.42857
$ws.Cells.Item(22, 9).Value = 242.42857  # I22: 242.57143 -> 242.42857
$ws.Cells.Item(22, 10).Value = 0  # J22: 120 -> 0
$ws.Cells.Item(22, 11).Value = 242.42857  # K22: 242.57143 -> 242.42857
$ws.Cells.Item(22, 12).Value = 0  # L22: 120 -> 0
$ws.Cells.Item(22, 13).Value = -69.42857000000001  # M22: -69.57142999999999 -> -69.42857000000001
$ws.Cells.Item(22, 14).ClearContents()  # N22: -466 -> (removed)
$ws.Cells.Item(80, 8).Value = 382.36365  # H80: 368.13043 -> 382.36365
$ws.Cells.Item(80, 10).Value = 436.625  # J80: 414.17648 -> 436.625
$ws.Cells.Item(80, 12).Value = 436.625  # L80: 414.17648 -> 436.625
$ws.Cells.Item(80, 14).Value = -2432.625  # N80: -2410.17648 -> -2432.625
$ws.Cells.Item(83, 8).Value = 382.36365  # H83: 368.13043 -> 382.36365
$ws.Cells.Item(83, 10).Value = 436.625  # J83: 414.17648 -> 436.625
$ws.Cells.Item(83, 12).Value = 2183.125  # L83: 2070.8824 -> 2183.125
$ws.Cells.Item(83, 14).Value = -12167.125  # N83: -12054.8824 -> -12167.125
$ws.Cells.Item(86, 8).Value = 8626.120000000001  # H86: 8691.32 -> 8626.120000000001
$ws.Cells.Item(86, 9).Value = 8715.166999999999  # I86: 8805.723 -> 8715.166999999999
$ws.Cells.Item(86, 11).Value = 8715.166999999999  # K86: 8805.723 -> 8715.166999999999
$ws.Cells.Item(86, 13).Value = -7592.166999999999  # M86: -7682.723 -> -7592.166999999999
$ws.Cells.Item(89, 8).Value = 8626.120000000001  # H89: 8691.32 -> 8626.120000000001
$ws.Cells.Item(89, 9).Value = 8715.166999999999  # I89: 8805.723 -> 8715.166999999999
$ws.Cells.Item(89, 11).Value = 43575.835  # K89: 44028.615 -> 43575.835
$ws.Cells.Item(89, 13).Value = -37959.835  # M89: -38412.615 -> -37959.835
$ws.Cells.Item(94, 8).Value = 803.7273  # H94: 867.1429000000001 -> 803.7273
$ws.Cells.Item(94, 9).Value = 803.7273  # I94: 860.15 -> 803.7273
$ws.Cells.Item(94, 10).Value = 0  # J94: 1007 -> 0
$ws.Cells.Item(94, 11).Value = 803.7273  # K94: 860.15 -> 803.7273
$ws.Cells.Item(94, 12).Value = 0  # L94: 1007 -> 0
$ws.Cells.Item(94, 13).Value = -352.7273  # M94: -409.15 -> -352.7273
$ws.Cells.Item(94, 14).ClearContents()  # N94: -1909 -> (removed)
$ws.Cells.Item(99, 8).Value = 1999.5  # H99: 892.5714 -> 1999.5
$ws.Cells.Item(99, 9).Value = 0  # I99: 449.6 -> 0
$ws.Cells.Item(99, 10).Value = 1999.5  # J99: 2000 -> 1999.5
$ws.Cells.Item(99, 11).Value = 0  # K99: 449.6 -> 0
$ws.Cells.Item(99, 12).Value = 1999.5  # L99: 2000 -> 1999.5
$ws.Cells.Item(99, 13).ClearContents()  # M99: 1048.4 -> (removed)
$ws.Cells.Item(99, 14).Value = -4995.5  # N99: -4996 -> -4995.5
$ws.Cells.Item(107, 8).Value = 3175.2273  # H107: 3142.8696 -> 3175.2273
$ws.Cells.Item(107, 9).Value = 3175.2273  # I107: 3184.7727 -> 3175.2273
$ws.Cells.Item(107, 10).Value = 0  # J107: 2221 -> 0
$ws.Cells.Item(107, 11).Value = 3175.2273  # K107: 3184.7727 -> 3175.2273
$ws.Cells.Item(107, 12).Value = 0  # L107: 2221 -> 0
$ws.Cells.Item(107, 13).Value = -1255.2273  # M107: -1264.7727 -> -1255.2273
$ws.Cells.Item(107, 14).ClearContents()  # N107: -6061 -> (removed)
$ws.Cells.Item(134, 8).Value = 4904.1665  # H134: 5016.2354 -> 4904.1665
$ws.Cells.Item(134, 10).Value = 4574.875  # J134: 4800 -> 4574.875
$ws.Cells.Item(134, 12).Value = 13724.625  # L134: 14400 -> 13724.625
$ws.Cells.Item(134, 14).Value = -18794.625  # N134: -19470 -> -18794.625
$ws.Cells.Item(137, 8).Value = 86666.336  # H137: 86666.664 -> 86666.336
$ws.Cells.Item(137, 10).Value = 86666.336  # J137: 86666.664 -> 86666.336
$ws.Cells.Item(137, 12).Value = 86666.336  # L137: 86666.664 -> 86666.336
$ws.Cells.Item(137, 14).Value = -96866.336  # N137: -96866.664 -> -96866.336

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 3466.4768  # H31: 3553.9834 -> 3466.4768
$ws.Cells.Item(31, 9).Value = 2346.7693  # I31: 2300.1667 -> 2346.7693
$ws.Cells.Item(31, 10).Value = 3746.4038  # J31: 3867.4375 -> 3746.4038
$ws.Cells.Item(31, 11).Value = 2346.7693  # K31: 2300.1667 -> 2346.7693
$ws.Cells.Item(31, 12).Value = 3746.4038  # L31: 3867.4375 -> 3746.4038
$ws.Cells.Item(31, 13).Value = -2051.7693  # M31: -2005.1667 -> -2051.7693
$ws.Cells.Item(31, 14).Value = -4336.4038  # N31: -4457.4375 -> -4336.4038
$ws.Cells.Item(34, 8).Value = 3466.4768  # H34: 3553.9834 -> 3466.4768
$ws.Cells.Item(34, 9).Value = 2346.7693  # I34: 2300.1667 -> 2346.7693
$ws.Cells.Item(34, 10).Value = 3746.4038  # J34: 3867.4375 -> 3746.4038
$ws.Cells.Item(34, 11).Value = 2346.7693  # K34: 2300.1667 -> 2346.7693
$ws.Cells.Item(34, 12).Value = 3746.4038  # L34: 3867.4375 -> 3746.4038
$ws.Cells.Item(34, 13).Value = -2144.7693  # M34: -2098.1667 -> -2144.7693
$ws.Cells.Item(34, 14).Value = -4150.4038  # N34: -4271.4375 -> -4150.4038
$ws.Cells.Item(58, 8).Value = 8352.049999999999  # H58: 8375.842000000001 -> 8352.049999999999
$ws.Cells.Item(58, 10).Value = 8642.5  # J58: 8710 -> 8642.5
$ws.Cells.Item(58, 12).Value = 8642.5  # L58: 8710 -> 8642.5
$ws.Cells.Item(58, 14).Value = -9048.5  # N58: -9116 -> -9048.5
$ws.Cells.Item(99, 8).Value = 5976  # H99: 3872.7742 -> 5976
$ws.Cells.Item(99, 9).Value = 5849.5  # I99: 3799.2068 -> 5849.5
$ws.Cells.Item(99, 10).Value = 8000  # J99: 4939.5 -> 8000
$ws.Cells.Item(99, 11).Value = 5849.5  # K99: 3799.2068 -> 5849.5
$ws.Cells.Item(99, 12).Value = 8000  # L99: 4939.5 -> 8000
$ws.Cells.Item(99, 13).Value = -4351.5  # M99: -2301.2068 -> -4351.5
$ws.Cells.Item(99, 14).Value = -10996  # N99: -7935.5 -> -10996
$ws.Cells.Item(107, 8).Value = 668.94446  # H107: 648.2222 -> 668.94446
$ws.Cells.Item(107, 9).Value = 502.91666  # I107: 469.27274 -> 502.91666
$ws.Cells.Item(107, 10).Value = 1001  # J107: 929.4286 -> 1001
$ws.Cells.Item(107, 11).Value = 502.91666  # K107: 469.27274 -> 502.91666
$ws.Cells.Item(107, 12).Value = 1001  # L107: 929.4286 -> 1001
$ws.Cells.Item(107, 13).Value = 1417.08334  # M107: 1450.72726 -> 1417.08334
$ws.Cells.Item(107, 14).Value = -4841  # N107: -4769.4286 -> -4841
$ws.Cells.Item(126, 8).Value = 5976  # H126: 3872.7742 -> 5976
$ws.Cells.Item(126, 9).Value = 5849.5  # I126: 3799.2068 -> 5849.5
$ws.Cells.Item(126, 10).Value = 8000  # J126: 4939.5 -> 8000
$ws.Cells.Item(126, 11).Value = 17548.5  # K126: 11397.6204 -> 17548.5
$ws.Cells.Item(126, 12).Value = 24000  # L126: 14818.5 -> 24000
$ws.Cells.Item(126, 13).Value = -15078.5  # M126: -8927.6204 -> -15078.5
$ws.Cells.Item(126, 14).Value = -28940  # N126: -19758.5 -> -28940
$ws.Cells.Item(132, 8).Value = 3774.25  # H132: 876.6667 -> 3774.25
$ws.Cells.Item(132, 9).Value = 1199.1666  # I132: 820 -> 1199.1666
$ws.Cells.Item(132, 10).Value = 11499.5  # J132: 1500 -> 11499.5
$ws.Cells.Item(132, 11).Value = 3597.4998  # K132: 2460 -> 3597.4998
$ws.Cells.Item(132, 12).Value = 34498.5  # L132: 4500 -> 34498.5
$ws.Cells.Item(132, 13).Value = -1067.4998  # M132: 70 -> -1067.4998
$ws.Cells.Item(132, 14).Value = -39558.5  # N132: -9560 -> -39558.5
$ws.Cells.Item(134, 8).Value = 6366.6875  # H134: 6279.1763 -> 6366.6875
$ws.Cells.Item(134, 9).Value = 6275.2144  # I134: 6182.1333 -> 6275.2144
$ws.Cells.Item(134, 11).Value = 18825.6432  # K134: 18546.3999 -> 18825.6432
$ws.Cells.Item(134, 13).Value = -16290.6432  # M134: -16011.3999 -> -16290.6432
$ws.Cells.Item(136, 8).Value = 8352.049999999999  # H136: 8375.842000000001 -> 8352.049999999999
$ws.Cells.Item(136, 10).Value = 8642.5  # J136: 8710 -> 8642.5
$ws.Cells.Item(136, 12).Value = 25927.5  # L136: 26130 -> 25927.5
$ws.Cells.Item(136, 14).Value = -31027.5  # N136: -31230 -> -31027.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(4, 8).Value = 44217750  # H4: 51454050 -> 44217750
$ws.Cells.Item(4, 9).Value = 56526350  # I4: 70457950 -> 56526350
$ws.Cells.Item(4, 11).Value = 169579050  # K4: 211373850 -> 169579050
$ws.Cells.Item(4, 13).Value = -169578938  # M4: -211373738 -> -169578938
$ws.Cells.Item(31, 8).Value = 5098  # H31: 5797.4 -> 5098
$ws.Cells.Item(31, 9).Value = 4397.3335  # I31: 4398.5 -> 4397.3335
$ws.Cells.Item(31, 10).Value = 6149  # J31: 6730 -> 6149
$ws.Cells.Item(31, 11).Value = 13192.0005  # K31: 13195.5 -> 13192.0005
$ws.Cells.Item(31, 12).Value = 18447  # L31: 20190 -> 18447
$ws.Cells.Item(31, 13).Value = -12904.0005  # M31: -12907.5 -> -12904.0005
$ws.Cells.Item(31, 14).Value = -19023  # N31: -20766 -> -19023
$ws.Cells.Item(55, 8).Value = 4255.75  # H55: 4835.2856 -> 4255.75
$ws.Cells.Item(55, 10).Value = 10466.333  # J55: 15600 -> 10466.333
$ws.Cells.Item(55, 12).Value = 31398.999  # L55: 46800 -> 31398.999
$ws.Cells.Item(55, 14).Value = -31752.999  # N55: -47154 -> -31752.999
$ws.Cells.Item(64, 8).Value = 2650  # H64: 1500 -> 2650
$ws.Cells.Item(64, 10).Value = 3800  # J64: 0 -> 3800
$ws.Cells.Item(64, 12).Value = 11400  # L64: 0 -> 11400
$ws.Cells.Item(64, 14).Value = -11940  # N64: None -> -11940
$ws.Cells.Item(67, 8).Value = 2650  # H67: 1500 -> 2650
$ws.Cells.Item(67, 10).Value = 3800  # J67: 0 -> 3800
$ws.Cells.Item(67, 12).Value = 11400  # L67: 0 -> 11400
$ws.Cells.Item(67, 14).Value = -13272  # N67: None -> -13272
$ws.Cells.Item(92, 8).Value = 304.85715  # H92: 305.10526 -> 304.85715
$ws.Cells.Item(92, 9).Value = 266.8  # I92: 266.75 -> 266.8
$ws.Cells.Item(92, 10).Value = 400  # J92: 370.85715 -> 400
$ws.Cells.Item(92, 11).Value = 800.4000000000001  # K92: 800.25 -> 800.4000000000001
$ws.Cells.Item(92, 12).Value = 1200  # L92: 1112.57145 -> 1200
$ws.Cells.Item(92, 13).Value = 447.5999999999999  # M92: 447.75 -> 447.5999999999999
$ws.Cells.Item(92, 14).Value = -3696  # N92: -3608.57145 -> -3696
$ws.Cells.Item(106, 8).Value = 15049.223  # H106: 0 -> 15049.223
$ws.Cells.Item(106, 10).Value = 15049.223  # J106: 0 -> 15049.223
$ws.Cells.Item(106, 12).Value = 45147.669  # L106: 0 -> 45147.669
$ws.Cells.Item(106, 14).Value = -47039.669  # N106: None -> -47039.669
$ws.Cells.Item(108, 8).Value = 422  # H108: 681.375 -> 422
$ws.Cells.Item(108, 9).Value = 422  # I108: 681.375 -> 422
$ws.Cells.Item(108, 11).Value = 1266  # K108: 2044.125 -> 1266
$ws.Cells.Item(108, 13).Value = 1614  # M108: 835.875 -> 1614
$ws.Cells.Item(110, 8).Value = 278931  # H110: 223290.2 -> 278931
$ws.Cells.Item(110, 9).Value = 7862  # I110: 5483.6665 -> 7862
$ws.Cells.Item(110, 11).Value = 23586  # K110: 16450.9995 -> 23586
$ws.Cells.Item(110, 13).Value = -19496  # M110: -12360.9995 -> -19496
$ws.Cells.Item(129, 8).Value = 2548.647  # H129: 2667 -> 2548.647
$ws.Cells.Item(129, 9).Value = 419.44446  # I129: 441.25 -> 419.44446
$ws.Cells.Item(129, 10).Value = 4944  # J129: 4761.8237 -> 4944
$ws.Cells.Item(129, 11).Value = 1258.33338  # K129: 1323.75 -> 1258.33338
$ws.Cells.Item(129, 12).Value = 14832  # L129: 14285.4711 -> 14832
$ws.Cells.Item(129, 13).Value = 3741.66662  # M129: 3676.25 -> 3741.66662
$ws.Cells.Item(129, 14).Value = -24832  # N129: -24285.4711 -> -24832
$ws.Cells.Item(137, 8).Value = 20840828  # H137: 17864034 -> 20840828
$ws.Cells.Item(137, 9).Value = 41670508  # I137: 35717790 -> 41670508
$ws.Cells.Item(137, 10).Value = 11148  # J137: 10274.429 -> 11148
$ws.Cells.Item(137, 11).Value = 125011524  # K137: 107153370 -> 125011524
$ws.Cells.Item(137, 12).Value = 33444  # L137: 30823.287 -> 33444
$ws.Cells.Item(137, 13).Value = -125006424  # M137: -107148270 -> -125006424
$ws.Cells.Item(137, 14).Value = -43644  # N137: -41023.287 -> -43644

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(24, 8).Value = 10250.5  # H24: 9712 -> 10250.5
$ws.Cells.Item(24, 10).Value = 10400  # J24: 9495 -> 10400
$ws.Cells.Item(24, 12).Value = 10400  # L24: 9495 -> 10400
$ws.Cells.Item(24, 14).Value = -10746  # N24: -9841 -> -10746
$ws.Cells.Item(70, 8).Value = 7664.467  # H70: 7926.9287 -> 7664.467
$ws.Cells.Item(70, 9).Value = 4497.6665  # I70: 4599.2 -> 4497.6665
$ws.Cells.Item(70, 11).Value = 4497.6665  # K70: 4599.2 -> 4497.6665
$ws.Cells.Item(70, 13).Value = -4227.6665  # M70: -4329.2 -> -4227.6665
$ws.Cells.Item(73, 8).Value = 7664.467  # H73: 7926.9287 -> 7664.467
$ws.Cells.Item(73, 9).Value = 4497.6665  # I73: 4599.2 -> 4497.6665
$ws.Cells.Item(73, 11).Value = 4497.6665  # K73: 4599.2 -> 4497.6665
$ws.Cells.Item(73, 13).Value = -3561.6665  # M73: -3663.2 -> -3561.6665
$ws.Cells.Item(92, 8).Value = 19949.818  # H92: 19370.75 -> 19949.818
$ws.Cells.Item(92, 10).Value = 19949.818  # J92: 19370.75 -> 19949.818
$ws.Cells.Item(92, 12).Value = 19949.818  # L92: 19370.75 -> 19949.818
$ws.Cells.Item(92, 14).Value = -23693.818  # N92: -23114.75 -> -23693.818
$ws.Cells.Item(102, 8).Value = 46175.918  # H102: 37758.734 -> 46175.918
$ws.Cells.Item(102, 9).Value = 4428.2856  # I102: 4326.8 -> 4428.2856
$ws.Cells.Item(102, 11).Value = 4428.2856  # K102: 4326.8 -> 4428.2856
$ws.Cells.Item(102, 13).Value = -2806.2856  # M102: -2704.8 -> -2806.2856
$ws.Cells.Item(122, 8).Value = 38540956  # H122: 40082510 -> 38540956
$ws.Cells.Item(122, 10).Value = 3000  # J122: 3500 -> 3000
$ws.Cells.Item(122, 12).Value = 9000  # L122: 10500 -> 9000
$ws.Cells.Item(122, 14).Value = -13900  # N122: -15400 -> -13900
$ws.Cells.Item(126, 8).Value = 5063.222  # H126: 5084.3887 -> 5063.222
$ws.Cells.Item(126, 9).Value = 3341.9167  # I126: 3358.5833 -> 3341.9167
$ws.Cells.Item(126, 10).Value = 8505.833000000001  # J126: 8536 -> 8505.833000000001
$ws.Cells.Item(126, 11).Value = 10025.7501  # K126: 10075.7499 -> 10025.7501
$ws.Cells.Item(126, 12).Value = 25517.499  # L126: 25608 -> 25517.499
$ws.Cells.Item(126, 13).Value = -7555.750100000001  # M126: -7605.749899999999 -> -7555.750100000001
$ws.Cells.Item(126, 14).Value = -30457.499  # N126: -30548 -> -30457.499
$ws.Cells.Item(132, 8).Value = 6989.65  # H132: 8349.625 -> 6989.65
$ws.Cells.Item(132, 9).Value = 9223.076999999999  # I132: 12633.444 -> 9223.076999999999
$ws.Cells.Item(132, 11).Value = 27669.231  # K132: 37900.33199999999 -> 27669.231
$ws.Cells.Item(132, 13).Value = -25139.231  # M132: -35370.33199999999 -> -25139.231

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(7, 8).Value = 6458251.5  # H7: 6458873 -> 6458251.5
$ws.Cells.Item(7, 9).Value = 8338774  # I7: 8005435 -> 8338774
$ws.Cells.Item(7, 10).Value = 10746.286  # J7: 14864.833 -> 10746.286
$ws.Cells.Item(7, 11).Value = 8338774  # K7: 8005435 -> 8338774
$ws.Cells.Item(7, 12).Value = 10746.286  # L7: 14864.833 -> 10746.286
$ws.Cells.Item(7, 13).Value = -8338662  # M7: -8005323 -> -8338662
$ws.Cells.Item(7, 14).Value = -10970.286  # N7: -15088.833 -> -10970.286
$ws.Cells.Item(39, 8).Value = 10000  # H39: 0 -> 10000
$ws.Cells.Item(39, 9).Value = 10000  # I39: 0 -> 10000
$ws.Cells.Item(39, 11).Value = 10000  # K39: 0 -> 10000
$ws.Cells.Item(39, 13).Value = -9540  # M39: None -> -9540
$ws.Cells.Item(40, 8).Value = 7660  # H40: 6180 -> 7660
$ws.Cells.Item(40, 9).Value = 10495  # I40: 7227.5 -> 10495
$ws.Cells.Item(40, 11).Value = 10495  # K40: 7227.5 -> 10495
$ws.Cells.Item(40, 13).Value = -10359  # M40: -7091.5 -> -10359
$ws.Cells.Item(46, 8).Value = 52084.32  # H46: 105650.086 -> 52084.32
$ws.Cells.Item(46, 9).Value = 5863  # I46: 19750 -> 5863
$ws.Cells.Item(46, 10).Value = 73835.53  # J46: 122830.1 -> 73835.53
$ws.Cells.Item(46, 11).Value = 5863  # K46: 19750 -> 5863
$ws.Cells.Item(46, 12).Value = 73835.53  # L46: 122830.1 -> 73835.53
$ws.Cells.Item(46, 13).Value = -5675  # M46: -19562 -> -5675
$ws.Cells.Item(46, 14).Value = -74211.53  # N46: -123206.1 -> -74211.53
$ws.Cells.Item(61, 8).Value = 9410.190000000001  # H61: 9636 -> 9410.190000000001
$ws.Cells.Item(61, 9).Value = 8655.4375  # I61: 8951.8125 -> 8655.4375
$ws.Cells.Item(61, 11).Value = 8655.4375  # K61: 8951.8125 -> 8655.4375
$ws.Cells.Item(61, 13).Value = -8453.4375  # M61: -8749.8125 -> -8453.4375
$ws.Cells.Item(93, 8).Value = 620.2692  # H93: 673 -> 620.2692
$ws.Cells.Item(93, 9).Value = 660.5625  # I93: 726.4286 -> 660.5625
$ws.Cells.Item(93, 10).Value = 555.8  # J93: 589.8889 -> 555.8
$ws.Cells.Item(93, 11).Value = 660.5625  # K93: 726.4286 -> 660.5625
$ws.Cells.Item(93, 12).Value = 555.8  # L93: 589.8889 -> 555.8
$ws.Cells.Item(93, 13).Value = 587.4375  # M93: 521.5714 -> 587.4375
$ws.Cells.Item(93, 14).Value = -3051.8  # N93: -3085.8889 -> -3051.8
$ws.Cells.Item(100, 8).Value = 3696.0715  # H100: 3782.9167 -> 3696.0715
$ws.Cells.Item(100, 9).Value = 3406.25  # I100: 3483.3333 -> 3406.25
$ws.Cells.Item(100, 11).Value = 3406.25  # K100: 3483.3333 -> 3406.25
$ws.Cells.Item(100, 13).Value = -2865.25  # M100: -2942.3333 -> -2865.25
$ws.Cells.Item(113, 8).Value = 9410.190000000001  # H113: 9636 -> 9410.190000000001
$ws.Cells.Item(113, 9).Value = 8655.4375  # I113: 8951.8125 -> 8655.4375
$ws.Cells.Item(113, 11).Value = 8655.4375  # K113: 8951.8125 -> 8655.4375
$ws.Cells.Item(113, 13).Value = -6485.4375  # M113: -6781.8125 -> -6485.4375
$ws.Cells.Item(122, 8).Value = 8958.450000000001  # H122: 9232.842000000001 -> 8958.450000000001
$ws.Cells.Item(122, 9).Value = 7388.727  # I122: 7753.1 -> 7388.727
$ws.Cells.Item(122, 11).Value = 22166.181  # K122: 23259.3 -> 22166.181
$ws.Cells.Item(122, 13).Value = -19716.181  # M122: -20809.3 -> -19716.181
$ws.Cells.Item(126, 8).Value = 6458251.5  # H126: 6458873 -> 6458251.5
$ws.Cells.Item(126, 9).Value = 8338774  # I126: 8005435 -> 8338774
$ws.Cells.Item(126, 10).Value = 10746.286  # J126: 14864.833 -> 10746.286
$ws.Cells.Item(126, 11).Value = 25016322  # K126: 24016305 -> 25016322
$ws.Cells.Item(126, 12).Value = 32238.858  # L126: 44594.499 -> 32238.858
$ws.Cells.Item(126, 13).Value = -25013852  # M126: -24013835 -> -25013852
$ws.Cells.Item(126, 14).Value = -37178.858  # N126: -49534.499 -> -37178.858
$ws.Cells.Item(132, 8).Value = 7339.3257  # H132: 7497.405 -> 7339.3257
$ws.Cells.Item(132, 9).Value = 7017.3237  # I132: 7208.758 -> 7017.3237
$ws.Cells.Item(132, 11).Value = 21051.9711  # K132: 21626.274 -> 21051.9711
$ws.Cells.Item(132, 13).Value = -18521.9711  # M132: -19096.274 -> -18521.9711
$ws.Cells.Item(133, 8).Value = 103750  # H133: 149499 -> 103750
$ws.Cells.Item(133, 10).Value = 103750  # J133: 149499 -> 103750
$ws.Cells.Item(133, 12).Value = 103750  # L133: 149499 -> 103750
$ws.Cells.Item(133, 14).Value = -108810  # N133: -154559 -> -108810
$ws.Cells.Item(136, 8).Value = 3864.3784  # H136: 4059.7942 -> 3864.3784
$ws.Cells.Item(136, 9).Value = 3504.7144  # I136: 3725.32 -> 3504.7144
$ws.Cells.Item(136, 10).Value = 4983.3335  # J136: 4988.8887 -> 4983.3335
$ws.Cells.Item(136, 11).Value = 10514.1432  # K136: 11175.96 -> 10514.1432
$ws.Cells.Item(136, 12).Value = 14950.0005  # L136: 14966.6661 -> 14950.0005
$ws.Cells.Item(136, 13).Value = -7964.143199999999  # M136: -8625.960000000001 -> -7964.143199999999
$ws.Cells.Item(136, 14).Value = -20050.0005  # N136: -20066.6661 -> -20050.0005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(31, 8).Value = 0  # H31: 12500 -> 0
$ws.Cells.Item(31, 9).Value = 0  # I31: 5000 -> 0
$ws.Cells.Item(31, 10).Value = 0  # J31: 20000 -> 0
$ws.Cells.Item(31, 11).Value = 0  # K31: 5000 -> 0
$ws.Cells.Item(31, 12).Value = 0  # L31: 20000 -> 0
$ws.Cells.Item(31, 13).ClearContents()  # M31: -4652 -> (removed)
$ws.Cells.Item(31, 14).ClearContents()  # N31: -20696 -> (removed)
$ws.Cells.Item(45, 8).Value = 23597.4  # H45: 23657.2 -> 23597.4
$ws.Cells.Item(45, 10).Value = 27430.666  # J45: 27530.334 -> 27430.666
$ws.Cells.Item(45, 12).Value = 27430.666  # L45: 27530.334 -> 27430.666
$ws.Cells.Item(45, 14).Value = -28412.666  # N45: -28512.334 -> -28412.666
$ws.Cells.Item(54, 8).Value = 0  # H54: 50000 -> 0
$ws.Cells.Item(54, 10).Value = 0  # J54: 50000 -> 0
$ws.Cells.Item(54, 12).Value = 0  # L54: 50000 -> 0
$ws.Cells.Item(54, 14).ClearContents()  # N54: -51040 -> (removed)
$ws.Cells.Item(62, 8).Value = 75706.37  # H62: 60461.285 -> 75706.37
$ws.Cells.Item(62, 9).Value = 133677.83  # I62: 101232 -> 133677.83
$ws.Cells.Item(62, 10).Value = 6140.6  # J62: 6100.3335 -> 6140.6
$ws.Cells.Item(62, 11).Value = 133677.83  # K62: 101232 -> 133677.83
$ws.Cells.Item(62, 12).Value = 6140.6  # L62: 6100.3335 -> 6140.6
$ws.Cells.Item(62, 13).Value = -133053.83  # M62: -100608 -> -133053.83
$ws.Cells.Item(62, 14).Value = -7388.6  # N62: -7348.3335 -> -7388.6
$ws.Cells.Item(65, 8).Value = 75706.37  # H65: 60461.285 -> 75706.37
$ws.Cells.Item(65, 9).Value = 133677.83  # I65: 101232 -> 133677.83
$ws.Cells.Item(65, 10).Value = 6140.6  # J65: 6100.3335 -> 6140.6
$ws.Cells.Item(65, 11).Value = 668389.1499999999  # K65: 506160 -> 668389.1499999999
$ws.Cells.Item(65, 12).Value = 30703  # L65: 30501.6675 -> 30703
$ws.Cells.Item(65, 13).Value = -665269.1499999999  # M65: -503040 -> -665269.1499999999
$ws.Cells.Item(65, 14).Value = -36943  # N65: -36741.6675 -> -36943
$ws.Cells.Item(81, 8).Value = 6926.9165  # H81: 7417.5 -> 6926.9165
$ws.Cells.Item(81, 9).Value = 5060  # I81: 5277.25 -> 5060
$ws.Cells.Item(81, 10).Value = 7860.375  # J81: 8640.5 -> 7860.375
$ws.Cells.Item(81, 11).Value = 10120  # K81: 10554.5 -> 10120
$ws.Cells.Item(81, 12).Value = 15720.75  # L81: 17281 -> 15720.75
$ws.Cells.Item(81, 13).Value = -9059  # M81: -9493.5 -> -9059
$ws.Cells.Item(81, 14).Value = -17842.75  # N81: -19403 -> -17842.75
$ws.Cells.Item(84, 8).Value = 6926.9165  # H84: 7417.5 -> 6926.9165
$ws.Cells.Item(84, 9).Value = 5060  # I84: 5277.25 -> 5060
$ws.Cells.Item(84, 10).Value = 7860.375  # J84: 8640.5 -> 7860.375
$ws.Cells.Item(84, 11).Value = 50600  # K84: 52772.5 -> 50600
$ws.Cells.Item(84, 12).Value = 78603.75  # L84: 86405 -> 78603.75
$ws.Cells.Item(84, 13).Value = -45296  # M84: -47468.5 -> -45296
$ws.Cells.Item(84, 14).Value = -89211.75  # N84: -97013 -> -89211.75
$ws.Cells.Item(107, 8).Value = 731.26086  # H107: 762.43475 -> 731.26086
$ws.Cells.Item(107, 9).Value = 594.9375  # I107: 672.53845 -> 594.9375
$ws.Cells.Item(107, 10).Value = 1042.8572  # J107: 879.3 -> 1042.8572
$ws.Cells.Item(107, 11).Value = 1784.8125  # K107: 2017.61535 -> 1784.8125
$ws.Cells.Item(107, 12).Value = 3128.5716  # L107: 2637.9 -> 3128.5716
$ws.Cells.Item(107, 13).Value = 135.1875  # M107: -97.61535000000003 -> 135.1875
$ws.Cells.Item(107, 14).Value = -6968.571599999999  # N107: -6477.9 -> -6968.571599999999
$ws.Cells.Item(109, 8).Value = 0  # H109: 96420 -> 0
$ws.Cells.Item(109, 10).Value = 0  # J109: 96420 -> 0
$ws.Cells.Item(109, 12).Value = 0  # L109: 96420 -> 0
$ws.Cells.Item(109, 14).ClearContents()  # N109: -99194 -> (removed)
$ws.Cells.Item(113, 8).Value = 5209153  # H113: 5209154.5 -> 5209153
$ws.Cells.Item(113, 10).Value = 1343.6  # J113: 1348 -> 1343.6
$ws.Cells.Item(113, 12).Value = 4030.8  # L113: 4044 -> 4030.8
$ws.Cells.Item(113, 14).Value = -8370.799999999999  # N113: -8384 -> -8370.799999999999
$ws.Cells.Item(122, 8).Value = 3374.5532  # H122: 3349.7346 -> 3374.5532
$ws.Cells.Item(122, 9).Value = 1925.2941  # I122: 1972.0278 -> 1925.2941
$ws.Cells.Item(122, 11).Value = 5775.8823  # K122: 5916.0834 -> 5775.8823
$ws.Cells.Item(122, 13).Value = -3325.8823  # M122: -3466.0834 -> -3325.8823
$ws.Cells.Item(126, 8).Value = 3237.3125  # H126: 3214.9697 -> 3237.3125
$ws.Cells.Item(126, 9).Value = 3300.3809  # I126: 3263.9092 -> 3300.3809
$ws.Cells.Item(126, 10).Value = 3116.9092  # J126: 3117.0908 -> 3116.9092
$ws.Cells.Item(126, 11).Value = 9901.1427  # K126: 9791.7276 -> 9901.1427
$ws.Cells.Item(126, 12).Value = 9350.7276  # L126: 9351.2724 -> 9350.7276
$ws.Cells.Item(126, 13).Value = -7431.1427  # M126: -7321.7276 -> -7431.1427
$ws.Cells.Item(126, 14).Value = -14290.7276  # N126: -14291.2724 -> -14290.7276
$ws.Cells.Item(132, 8).Value = 3573.4614  # H132: 3967.238 -> 3573.4614
$ws.Cells.Item(132, 9).Value = 3194.5386  # I132: 3991.375 -> 3194.5386
$ws.Cells.Item(132, 11).Value = 9583.6158  # K132: 11974.125 -> 9583.6158
$ws.Cells.Item(132, 13).Value = -7053.6158  # M132: -9444.125 -> -7053.6158
$ws.Cells.Item(136, 8).Value = 3718.1943  # H136: 3844.4285 -> 3718.1943
$ws.Cells.Item(136, 9).Value = 1968.3334  # I136: 2001.7391 -> 1968.3334
$ws.Cells.Item(136, 10).Value = 7217.9165  # J136: 7376.25 -> 7217.9165
$ws.Cells.Item(136, 11).Value = 5905.0002  # K136: 6005.2173 -> 5905.0002
$ws.Cells.Item(136, 12).Value = 21653.7495  # L136: 22128.75 -> 21653.7495
$ws.Cells.Item(136, 13).Value = -3355.0002  # M136: -3455.2173 -> -3355.0002
$ws.Cells.Item(136, 14).Value = -26753.7495  # N136: -27228.75 -> -26753.7495
